$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Icam2", "Itgal", "ECs", 3, 1, 42.675953, 128.027859, 0.9123907898314253, 0.9123907898314252, 1, 0.3333333333333333, 0.1050723333333333, 0.315217, 0.006764366142259137, 0.006764366142259137, 4.484061958933667, 40.356557630403, 0.006171745367244765, 0.006171745367244764),
  @("ECs", "Icam2", "Itgal", "FAPs", 3, 1, 42.675953, 128.027859, 0.9123907898314253, 0.9123907898314252, 3, 1, 0.9976189999999999, 2.992857, 0.06422490081252995, 0.06422490081252995, 42.574341555907, 383.169074003163, 0.05859820797918915, 0.05859820797918914),
  @("ECs", "Icam2", "Itgal", "MuSCs", 3, 1, 42.675953, 128.027859, 0.9123907898314253, 0.9123907898314252, 3, 1, 0.173444, 0.520332, 0.01116600996625811, 0.01116600996625811, 7.401887992132001, 66.616991929188, 0.0101877646523798, 0.0101877646523798),
  @("ECs", "Icam2", "Itgal", "Resolving-Mac", 3, 1, 42.675953, 128.027859, 0.9123907898314253, 0.9123907898314252, 3, 1, 14.25707666666667, 42.77123, 0.9178447230789528, 0.9178447230789528, 608.4343337440635, 5475.909003696571, 0.8374330718326116, 0.8374330718326115),
  @("FAPs", "Icam2", "Itgal", "ECs", 3, 1, 1.529744666666667, 4.589234, 0.03270518515803058, 0.03270518515803057, 1, 0.3333333333333333, 0.1050723333333333, 0.315217, 0.006764366142259137, 0.006764366142259137, 0.1607338415308889, 1.446604573778, 0.0002212298471592981, 0.000221229847159298),
  @("FAPs", "Icam2", "Itgal", "FAPs", 3, 1, 1.529744666666667, 4.589234, 0.03270518515803058, 0.03270518515803057, 3, 1, 0.9976189999999999, 2.992857, 0.06422490081252995, 0.06422490081252995, 1.526102344615333, 13.734921101538, 0.00210048727282994, 0.00210048727282994),
  @("FAPs", "Icam2", "Itgal", "MuSCs", 3, 1, 1.529744666666667, 4.589234, 0.03270518515803058, 0.03270518515803057, 3, 1, 0.173444, 0.520332, 0.01116600996625811, 0.01116600996625811, 0.2653250339653334, 2.387925305688, 0.0003651864234228861, 0.000365186423422886),
  @("FAPs", "Icam2", "Itgal", "Resolving-Mac", 3, 1, 1.529744666666667, 4.589234, 0.03270518515803058, 0.03270518515803057, 3, 1, 14.25707666666667, 42.77123, 0.9178447230789528, 0.9178447230789528, 21.80968699309112, 196.28718293782, 0.03001828161461845, 0.03001828161461845),
  @("MuSCs", "Icam2", "Itgal", "ECs", 3, 1, 0.9319226666666666, 2.795768, 0.01992404616955614, 0.01992404616955614, 1, 0.3333333333333333, 0.1050723333333333, 0.315217, 0.006764366142259137, 0.006764366142259137, 0.09791928907288888, 0.881273601656, 0.0001347735433261534, 0.0001347735433261534),
  @("MuSCs", "Icam2", "Itgal", "FAPs", 3, 1, 0.9319226666666666, 2.795768, 0.01992404616955614, 0.01992404616955614, 3, 1, 0.9976189999999999, 2.992857, 0.06422490081252995, 0.06422490081252995, 0.9297037587973331, 8.367333829176, 0.001279619889024011, 0.00127961988902401),
  @("MuSCs", "Icam2", "Itgal", "MuSCs", 3, 1, 0.9319226666666666, 2.795768, 0.01992404616955614, 0.01992404616955614, 3, 1, 0.173444, 0.520332, 0.01116600996625811, 0.01116600996625811, 0.1616363949973333, 1.454727554976, 0.0002224720980974506, 0.0002224720980974505),
  @("MuSCs", "Icam2", "Itgal", "Resolving-Mac", 3, 1, 0.9319226666666666, 2.795768, 0.01992404616955614, 0.01992404616955614, 3, 1, 14.25707666666667, 42.77123, 0.9178447230789528, 0.9178447230789528, 13.28649290607111, 119.57843615464, 0.01828718063910853, 0.01828718063910853),
  @("Resolving-Mac", "Icam2", "Itgal", "ECs", 3, 1, 1.636145333333333, 4.908436, 0.03497997884098805, 0.03497997884098804, 1, 0.3333333333333333, 0.1050723333333333, 0.315217, 0.006764366142259137, 0.006764366142259137, 0.1719136078457778, 1.547222470612, 0.0002366173845289206, 0.0002366173845289205),
  @("Resolving-Mac", "Icam2", "Itgal", "FAPs", 3, 1, 1.636145333333333, 4.908436, 0.03497997884098805, 0.03497997884098804, 3, 1, 0.9976189999999999, 2.992857, 0.06422490081252995, 0.06422490081252995, 1.632249671294667, 14.690247041652, 0.002246585671486854, 0.002246585671486853),
  @("Resolving-Mac", "Icam2", "Itgal", "MuSCs", 3, 1, 1.636145333333333, 4.908436, 0.03497997884098805, 0.03497997884098804, 3, 1, 0.173444, 0.520332, 0.01116600996625811, 0.01116600996625811, 0.2837795911946667, 2.554016320752, 0.0003905867923579703, 0.0003905867923579702),
  @("Resolving-Mac", "Icam2", "Itgal", "Resolving-Mac", 3, 1, 1.636145333333333, 4.908436, 0.03497997884098805, 0.03497997884098804, 3, 1, 14.25707666666667, 42.77123, 0.9178447230789528, 0.9178447230789528, 23.32664945514222, 209.93984509628, 0.0321061889926143, 0.0321061889926143)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  $r = $startRow + $i
  for ($j = 0; $j -lt $row.Length; $j++) {
    $ws.Cells.Item($r, $j + 1).Value = $row[$j]
  }
}
